# Insert a new data row at row 194 (shifts existing rows 194.. down by one)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("194:194").Insert()

# Populate the newly inserted row 194 with the new record
$ws.Range("A194").Value = 11
$ws.Range("B194").Value = "Vega Monumental Concepción"
$ws.Range("C194").Value = "Bíobío"
$ws.Range("D194").Value = 44518
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 100112020
$ws.Range("G194").Value = "Tomate"
$ws.Range("H194").Value = "Larga vida"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 500
$ws.Range("K194").Value = 8000
$ws.Range("L194").Value = 8500
$ws.Range("M194").Value = 8250
$ws.Range("N194").Value = "$/caja 12 kilos"
$ws.Range("O194").Value = "Región de Arica y Parinacota"
$ws.Range("P194").Value = 688
$ws.Range("Q194").Value = 12
$ws.Range("R194").Value = "Hortaliza"
